$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 3.2
$ws.Range("G2").Value = 3.25
$ws.Range("H2").Value = 2.42
$ws.Range("I2").Value = 2.44
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 1.9
$ws.Range("Q2").Value = 1.98
$ws.Range("S2").Value = 3.5
$ws.Range("T2").Value = 1.77
$ws.Range("U2").Value = 2.18
$ws.Range("V2").Value = 1.69
$ws.Range("W2").Value = 1.44
$ws.Range("Y2").Value = 10.5
$ws.Range("AC2").Value = 8
$ws.Range("AH2").Value = 17.5

# Row 3
$ws.Range("F3").Value = 1.84
$ws.Range("G3").Value = 2.04
$ws.Range("I3").Value = 4.8
$ws.Range("Q3").Value = 1.79

# Row 8
$ws.Range("F8").Value = 7.4
$ws.Range("G8").Value = 10.5
$ws.Range("H8").Value = 1.42
$ws.Range("I8").Value = 1.51
$ws.Range("J8").Value = 4.8
$ws.Range("K8").Value = 5.7
$ws.Range("P8").Value = 2.38
$ws.Range("Q8").Value = 1.59

# Row 9
$ws.Range("I9").Value = 2.18
$ws.Range("P9").Value = 2.34
$ws.Range("AB9").Value = 17.5
$ws.Range("AD9").Value = 11
$ws.Range("AH9").Value = 15.5

# Row 10
$ws.Range("AC10").Value = 6.8

$wb.Save()
